$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(11, 1).Value = "EGGW"
$ws.Cells.Item(11, 2).Value = -0.3683330118656158
$ws.Cells.Item(11, 3).Value = 51.87469863891602
$ws.Cells.Item(12, 1).Value = "EGKK"
$ws.Cells.Item(12, 2).Value = -0.190278
$ws.Cells.Item(12, 3).Value = 51.148102
$ws.Cells.Item(13, 1).Value = "EGLL"
$ws.Cells.Item(13, 2).Value = -0.461941
$ws.Cells.Item(13, 3).Value = 51.4706
$ws.Cells.Item(14, 1).Value = "EGSS"
$ws.Cells.Item(14, 2).Value = 0.234999999404
$ws.Cells.Item(14, 3).Value = 51.8849983215
$ws.Cells.Item(15, 1).Value = "EHAM"
$ws.Cells.Item(15, 2).Value = 4.76389
$ws.Cells.Item(15, 3).Value = 52.308601
$ws.Cells.Item(16, 1).Value = "EHEH"
$ws.Cells.Item(16, 2).Value = 5.374529838559999
$ws.Cells.Item(16, 3).Value = 51.4500999451
$ws.Cells.Item(17, 1).Value = "EINN"
$ws.Cells.Item(17, 2).Value = -8.92482
$ws.Cells.Item(17, 3).Value = 52.702
$ws.Cells.Item(18, 1).Value = "ELLX"
$ws.Cells.Item(18, 2).Value = 6.204444400000001
$ws.Cells.Item(18, 3).Value = 49.6233333
$ws.Cells.Item(19, 1).Value = "ENGM"
$ws.Cells.Item(19, 2).Value = 11.1004
$ws.Cells.Item(19, 3).Value = 60.193901
$ws.Cells.Item(20, 1).Value = "EVRA"
$ws.Cells.Item(20, 2).Value = 23.97109985351562
$ws.Cells.Item(20, 3).Value = 56.92359924316406
$ws.Cells.Item(21, 1).Value = "KORD"
$ws.Cells.Item(21, 2).Value = -87.90479999999999
$ws.Cells.Item(21, 3).Value = 41.9786
$ws.Cells.Item(22, 1).Value = "LCLK"
$ws.Cells.Item(22, 2).Value = 33.624901
$ws.Cells.Item(22, 3).Value = 34.875099
$ws.Cells.Item(23, 1).Value = "LDZA"
$ws.Cells.Item(23, 2).Value = 16.0687999725
$ws.Cells.Item(23, 3).Value = 45.74290084840001
$ws.Cells.Item(24, 1).Value = "LEMG"
$ws.Cells.Item(24, 2).Value = -4.49911
$ws.Cells.Item(24, 3).Value = 36.6749
$ws.Cells.Item(25, 1).Value = "LFMN"
$ws.Cells.Item(25, 2).Value = 7.215869903560001
$ws.Cells.Item(25, 3).Value = 43.6584014893
$ws.Cells.Item(26, 1).Value = "LFPG"
$ws.Cells.Item(26, 2).Value = 2.55
$ws.Cells.Item(26, 3).Value = 49.012798
$ws.Cells.Item(27, 1).Value = "LGAV"
$ws.Cells.Item(27, 2).Value = 23.9445
$ws.Cells.Item(27, 3).Value = 37.936401
$ws.Cells.Item(28, 1).Value = "LHBP"
$ws.Cells.Item(28, 2).Value = 19.261093
$ws.Cells.Item(28, 3).Value = 47.42976
$ws.Cells.Item(29, 1).Value = "LIMC"
$ws.Cells.Item(29, 2).Value = 8.728110000000001
$ws.Cells.Item(29, 3).Value = 45.6306
$ws.Cells.Item(30, 1).Value = "LIME"
$ws.Cells.Item(30, 2).Value = 9.70417
$ws.Cells.Item(30, 3).Value = 45.673901
$ws.Cells.Item(31, 1).Value = "LLBG"
$ws.Cells.Item(31, 2).Value = 34.88669967651367
$ws.Cells.Item(31, 3).Value = 32.01139831542969
$ws.Cells.Item(32, 1).Value = "LSGG"
$ws.Cells.Item(32, 2).Value = 6.108950138092041
$ws.Cells.Item(32, 3).Value = 46.23809814453125
$ws.Cells.Item(33, 1).Value = "LTAC"
$ws.Cells.Item(33, 2).Value = 33.011536
$ws.Cells.Item(33, 3).Value = 40.14624
$ws.Cells.Item(34, 1).Value = "LTAF"
$ws.Cells.Item(34, 2).Value = 35.280399
$ws.Cells.Item(34, 3).Value = 36.982201
$ws.Cells.Item(35, 1).Value = "LTAI"
$ws.Cells.Item(35, 2).Value = 30.800501
$ws.Cells.Item(35, 3).Value = 36.898701
$ws.Cells.Item(36, 1).Value = "LTAU"
$ws.Cells.Item(36, 2).Value = 35.494916
$ws.Cells.Item(36, 3).Value = 38.770248
$ws.Cells.Item(37, 1).Value = "LTBA"
$ws.Cells.Item(37, 2).Value = 28.823714
$ws.Cells.Item(37, 3).Value = 40.971913
$ws.Cells.Item(38, 1).Value = "LTBJ"
$ws.Cells.Item(38, 2).Value = 27.157
$ws.Cells.Item(38, 3).Value = 38.2924
$ws.Cells.Item(39, 1).Value = "LTBS"
$ws.Cells.Item(39, 2).Value = 28.7924995422
$ws.Cells.Item(39, 3).Value = 36.7131004333
$ws.Cells.Item(40, 1).Value = "LTFE"
$ws.Cells.Item(40, 2).Value = 27.6643009186
$ws.Cells.Item(40, 3).Value = 37.25059890749999
$ws.Cells.Item(41, 1).Value = "LTFG"
$ws.Cells.Item(41, 2).Value = 32.412707
$ws.Cells.Item(41, 3).Value = 36.208237
$ws.Cells.Item(42, 1).Value = "LTFJ"
$ws.Cells.Item(42, 2).Value = 29.3092
$ws.Cells.Item(42, 3).Value = 40.898602
$ws.Cells.Item(43, 1).Value = "LTFM"
$ws.Cells.Item(43, 2).Value = 28.741951
$ws.Cells.Item(43, 3).Value = 41.261297
$ws.Cells.Item(44, 1).Value = "LYBE"
$ws.Cells.Item(44, 2).Value = 20.3090991974
$ws.Cells.Item(44, 3).Value = 44.8184013367
$ws.Cells.Item(45, 1).Value = "LZIB"
$ws.Cells.Item(45, 2).Value = 17.21269989013672
$ws.Cells.Item(45, 3).Value = 48.17020034790039
$ws.Cells.Item(46, 1).Value = "OBBI"
$ws.Cells.Item(46, 2).Value = 50.63764
$ws.Cells.Item(46, 3).Value = 26.267295
$ws.Cells.Item(47, 1).Value = "OEDF"
$ws.Cells.Item(47, 2).Value = 49.79790115356445
$ws.Cells.Item(47, 3).Value = 26.47120094299316
$ws.Cells.Item(48, 1).Value = "OEJN"
$ws.Cells.Item(48, 2).Value = 39.156502
$ws.Cells.Item(48, 3).Value = 21.6796
$ws.Cells.Item(49, 1).Value = "OERK"
$ws.Cells.Item(49, 2).Value = 46.69879913330078
$ws.Cells.Item(49, 3).Value = 24.95759963989257
$ws.Cells.Item(50, 1).Value = "OIII"
$ws.Cells.Item(50, 2).Value = 51.322861
$ws.Cells.Item(50, 3).Value = 35.687622
$ws.Cells.Item(51, 1).Value = "OKBK"
$ws.Cells.Item(51, 2).Value = 47.96889877319336
$ws.Cells.Item(51, 3).Value = 29.22660064697266
$ws.Cells.Item(52, 1).Value = "OMAA"
$ws.Cells.Item(52, 2).Value = 54.651718
$ws.Cells.Item(52, 3).Value = 24.443764
$ws.Cells.Item(53, 1).Value = "OMDB"
$ws.Cells.Item(53, 2).Value = 55.3643989563
$ws.Cells.Item(53, 3).Value = 25.2527999878
$ws.Cells.Item(54, 1).Value = "OMDW"
$ws.Cells.Item(54, 2).Value = 55.161389
$ws.Cells.Item(54, 3).Value = 24.896356
$ws.Cells.Item(55, 1).Value = "OMSJ"
$ws.Cells.Item(55, 2).Value = 55.5172004699707
$ws.Cells.Item(55, 3).Value = 25.32859992980957
$ws.Cells.Item(56, 1).Value = "OPSR"
$ws.Cells.Item(56, 2).Value = 72.66500091552734
$ws.Cells.Item(56, 3).Value = 32.04859924316406
$ws.Cells.Item(57, 1).Value = "ORER"
$ws.Cells.Item(57, 2).Value = 43.96319961547852
$ws.Cells.Item(57, 3).Value = 36.23759841918945
$ws.Cells.Item(58, 1).Value = "OTHH"
$ws.Cells.Item(58, 2).Value = 51.608056
$ws.Cells.Item(58, 3).Value = 25.273056
$ws.Cells.Item(59, 1).Value = "RJAA"
$ws.Cells.Item(59, 2).Value = 140.386002
$ws.Cells.Item(59, 3).Value = 35.764702
$ws.Cells.Item(60, 1).Value = "RJBB"
$ws.Cells.Item(60, 2).Value = 135.244003
$ws.Cells.Item(60, 3).Value = 34.427299
$ws.Cells.Item(61, 1).Value = "RKSI"
$ws.Cells.Item(61, 2).Value = 126.4509963989258
$ws.Cells.Item(61, 3).Value = 37.46910095214844
$ws.Cells.Item(62, 1).Value = "UAAA"
$ws.Cells.Item(62, 2).Value = 77.06150700000001
$ws.Cells.Item(62, 3).Value = 43.364822
$ws.Cells.Item(63, 1).Value = "UBBB"
$ws.Cells.Item(63, 2).Value = 50.04669952392578
$ws.Cells.Item(63, 3).Value = 40.46749877929688
$ws.Cells.Item(64, 1).Value = "UCFM"
$ws.Cells.Item(64, 2).Value = 74.47760009770001
$ws.Cells.Item(64, 3).Value = 43.0612983704
$ws.Cells.Item(65, 1).Value = "UGTB"
$ws.Cells.Item(65, 2).Value = 44.95470047
$ws.Cells.Item(65, 3).Value = 41.6692008972
$ws.Cells.Item(66, 1).Value = "UKBB"
$ws.Cells.Item(66, 2).Value = 30.89469909667969
$ws.Cells.Item(66, 3).Value = 50.34500122070313
$ws.Cells.Item(67, 1).Value = "ULLI"
$ws.Cells.Item(67, 2).Value = 30.26250076293945
$ws.Cells.Item(67, 3).Value = 59.80030059814453
$ws.Cells.Item(68, 1).Value = "UMMS"
$ws.Cells.Item(68, 2).Value = 28.039964
$ws.Cells.Item(68, 3).Value = 53.888071
$ws.Cells.Item(69, 1).Value = "UNKL"
$ws.Cells.Item(69, 2).Value = 92.492437
$ws.Cells.Item(69, 3).Value = 56.173077
$ws.Cells.Item(70, 1).Value = "UNNT"
$ws.Cells.Item(70, 2).Value = 82.618675
$ws.Cells.Item(70, 3).Value = 55.01975600000001
$ws.Cells.Item(71, 1).Value = "URKK"
$ws.Cells.Item(71, 2).Value = 39.170501708984
$ws.Cells.Item(71, 3).Value = 45.034698486328
$ws.Cells.Item(72, 1).Value = "URMM"
$ws.Cells.Item(72, 2).Value = 43.08190155029297
$ws.Cells.Item(72, 3).Value = 44.22510147094727
$ws.Cells.Item(73, 1).Value = "URRP"
$ws.Cells.Item(73, 2).Value = 39.924722
$ws.Cells.Item(73, 3).Value = 47.493888
$ws.Cells.Item(74, 1).Value = "URWA"
$ws.Cells.Item(74, 2).Value = 48.0063018799
$ws.Cells.Item(74, 3).Value = 46.2832984924
$ws.Cells.Item(75, 1).Value = "URWW"
$ws.Cells.Item(75, 2).Value = 44.34550094604492
$ws.Cells.Item(75, 3).Value = 48.78250122070313
$ws.Cells.Item(76, 1).Value = "USPP"
$ws.Cells.Item(76, 2).Value = 56.021198
$ws.Cells.Item(76, 3).Value = 57.914501
$ws.Cells.Item(77, 1).Value = "USSS"
$ws.Cells.Item(77, 2).Value = 60.802700042725
$ws.Cells.Item(77, 3).Value = 56.743099212646
$ws.Cells.Item(78, 1).Value = "USTR"
$ws.Cells.Item(78, 2).Value = 65.3243026733
$ws.Cells.Item(78, 3).Value = 57.1896018982
$ws.Cells.Item(79, 1).Value = "UTSB"
$ws.Cells.Item(79, 2).Value = 64.483299
$ws.Cells.Item(79, 3).Value = 39.775002
$ws.Cells.Item(80, 1).Value = "UTTT"
$ws.Cells.Item(80, 2).Value = 69.258667
$ws.Cells.Item(80, 3).Value = 41.256088
$ws.Cells.Item(81, 1).Value = "UUBW"
$ws.Cells.Item(81, 2).Value = 38.150002
$ws.Cells.Item(81, 3).Value = 55.553299
$ws.Cells.Item(82, 1).Value = "UUDD"
$ws.Cells.Item(82, 2).Value = 37.90629959106445
$ws.Cells.Item(82, 3).Value = 55.40879821777344
$ws.Cells.Item(83, 1).Value = "UUEE"
$ws.Cells.Item(83, 2).Value = 37.4146
$ws.Cells.Item(83, 3).Value = 55.972599
$ws.Cells.Item(84, 1).Value = "UUOO"
$ws.Cells.Item(84, 2).Value = 39.22959899902344
$ws.Cells.Item(84, 3).Value = 51.81420135498047
$ws.Cells.Item(85, 1).Value = "UUWW"
$ws.Cells.Item(85, 2).Value = 37.2615013123
$ws.Cells.Item(85, 3).Value = 55.5914993286
$ws.Cells.Item(86, 1).Value = "UWGG"
$ws.Cells.Item(86, 2).Value = 43.784000396729
$ws.Cells.Item(86, 3).Value = 56.23009872436499
$ws.Cells.Item(87, 1).Value = "UWKD"
$ws.Cells.Item(87, 2).Value = 49.278701782227
$ws.Cells.Item(87, 3).Value = 55.606201171875
$ws.Cells.Item(88, 1).Value = "UWUU"
$ws.Cells.Item(88, 2).Value = 55.874401092529
$ws.Cells.Item(88, 3).Value = 54.55749893188501
$ws.Cells.Item(89, 1).Value = "UWWW"
$ws.Cells.Item(89, 2).Value = 50.16429901123
$ws.Cells.Item(89, 3).Value = 53.504901885986
$ws.Cells.Item(90, 1).Value = "VABB"
$ws.Cells.Item(90, 2).Value = 72.8678970337
$ws.Cells.Item(90, 3).Value = 19.0886993408
$ws.Cells.Item(91, 1).Value = "VGHS"
$ws.Cells.Item(91, 2).Value = 90.397783
$ws.Cells.Item(91, 3).Value = 23.843347
$ws.Cells.Item(92, 1).Value = "VHHH"
$ws.Cells.Item(92, 2).Value = 113.915001
$ws.Cells.Item(92, 3).Value = 22.308901
$ws.Cells.Item(93, 1).Value = "VIDP"
$ws.Cells.Item(93, 2).Value = 77.103104
$ws.Cells.Item(93, 3).Value = 28.5665
$ws.Cells.Item(94, 1).Value = "VVNB"
$ws.Cells.Item(94, 2).Value = 105.806999
$ws.Cells.Item(94, 3).Value = 21.221201
$ws.Cells.Item(95, 1).Value = "VVTS"
$ws.Cells.Item(95, 2).Value = 106.652
$ws.Cells.Item(95, 3).Value = 10.8188
$ws.Cells.Item(96, 1).Value = "WSSS"
$ws.Cells.Item(96, 2).Value = 103.994003
$ws.Cells.Item(96, 3).Value = 1.35019
$ws.Cells.Item(97, 1).Value = "ZBAA"
$ws.Cells.Item(97, 2).Value = 116.5849990844727
$ws.Cells.Item(97, 3).Value = 40.0801010131836
$ws.Cells.Item(98, 1).Value = "ZHCC"
$ws.Cells.Item(98, 2).Value = 113.849165
$ws.Cells.Item(98, 3).Value = 34.526497
$ws.Cells.Item(99, 1).Value = "ZSPD"
$ws.Cells.Item(99, 2).Value = 121.805
$ws.Cells.Item(99, 3).Value = 31.1434
